$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, shifting rows 99:163 down to 100:164
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new record.
$ws.Cells.Item(99, 1).Value = 9
$ws.Cells.Item(99, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(99, 3).Value = "Metropolitana"
$ws.Cells.Item(99, 4).Value = 45086
$ws.Cells.Item(99, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(99, 5).Value = 13
$ws.Cells.Item(99, 6).Value = 100112022
$ws.Cells.Item(99, 7).Value = "Arveja Verde"
$ws.Cells.Item(99, 8).Value = "Perfection"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 52
$ws.Cells.Item(99, 11).Value = 36000
$ws.Cells.Item(99, 12).Value = 38000
$ws.Cells.Item(99, 13).Value = 37000
$ws.Cells.Item(99, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(99, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(99, 16).Value = 1480
$ws.Cells.Item(99, 17).Value = 25
$ws.Cells.Item(99, 18).Value = "Hortaliza"
